$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from an existing date cell (A9) so the new
# date cells reuse the same cell style / number format instead of creating
# a duplicate numFmt entry.
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 10: date 25/09/2012 (serial 41177), Additional Effort [h] = 1
$ws.Cells.Item(10, 1).Value = 41177
$ws.Cells.Item(10, 3).Value = 1

# New row 11: date 26/09/2012 (serial 41178), Effort [h] = 2.5, comment text
$ws.Cells.Item(11, 1).Value = 41178
$ws.Cells.Item(11, 2).Value = 2.5
$ws.Cells.Item(11, 4).Value = "waitForEvent, setEvent implemented. TC03 added, but does not yet run"

$ws.Range("D11").Select()
